$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1038.25
$ws.Range("I29").Value = 701
$ws.Range("J29").Value = 2050
$ws.Range("K29").Value = 2103
$ws.Range("L29").Value = 6150
$ws.Range("M29").Value = -1822
$ws.Range("N29").Value = -6712

$ws.Range("H34").Value = 1507.6666
$ws.Range("I34").Value = 307.5
$ws.Range("J34").Value = 3908
$ws.Range("K34").Value = 307.5
$ws.Range("L34").Value = 3908
$ws.Range("M34").Value = -104.5
$ws.Range("N34").Value = -4314

$ws.Range("H36").Value = 1507.6666
$ws.Range("I36").Value = 307.5
$ws.Range("J36").Value = 3908
$ws.Range("K36").Value = 307.5
$ws.Range("L36").Value = 3908
$ws.Range("M36").Value = 407.5
$ws.Range("N36").Value = -5338

$ws.Range("H113").Value = 1690
$ws.Range("J113").Value = 1690
$ws.Range("L113").Value = 1690
$ws.Range("N113").Value = -8198

$ws.Range("H116").Value = 3653.4666
$ws.Range("I116").Value = 2525.25
$ws.Range("J116").Value = 4942.857
$ws.Range("K116").Value = 2525.25
$ws.Range("L116").Value = 4942.857
$ws.Range("M116").Value = 916.75
$ws.Range("N116").Value = -11826.857

$ws.Range("H129").Value = 957.12195
$ws.Range("I129").Value = 375
$ws.Range("J129").Value = 1037.9722
$ws.Range("K129").Value = 1125
$ws.Range("L129").Value = 3113.9166
$ws.Range("M129").Value = 3875
$ws.Range("N129").Value = -13113.9166

$ws.Range("H135").Value = 1552.725
$ws.Range("I135").Value = 349.75
$ws.Range("J135").Value = 4359.6665
$ws.Range("K135").Value = 3147.75
$ws.Range("L135").Value = 39236.9985
$ws.Range("M135").Value = -612.75
$ws.Range("N135").Value = -44306.9985

$ws.Range("H138").Value = 1615693.9
$ws.Range("I138").Value = 2859097.8
$ws.Range("J138").Value = 3873.926
$ws.Range("K138").Value = 8577293.399999999
$ws.Range("L138").Value = 11621.778
$ws.Range("M138").Value = -8572153.399999999
$ws.Range("N138").Value = -21901.778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3641.0667
$ws.Range("I16").Value = 4360
$ws.Range("J16").Value = 2203.2
$ws.Range("K16").Value = 4360
$ws.Range("L16").Value = 2203.2
$ws.Range("M16").Value = -4073
$ws.Range("N16").Value = -2777.2

$ws.Range("H113").Value = 3641.0667
$ws.Range("I113").Value = 4360
$ws.Range("J113").Value = 2203.2
$ws.Range("K113").Value = 4360
$ws.Range("L113").Value = 2203.2
$ws.Range("M113").Value = -2190
$ws.Range("N113").Value = -6543.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 1550
$ws.Range("J31").Value = 3300
$ws.Range("L31").Value = 9900
$ws.Range("N31").Value = -10476

$ws.Range("H35").Value = 2130
$ws.Range("J35").Value = 2130
$ws.Range("L35").Value = 6390
$ws.Range("N35").Value = -6966

$ws.Range("H97").Value = 3386.5715
$ws.Range("I97").Value = 153
$ws.Range("J97").Value = 4680
$ws.Range("K97").Value = 459
$ws.Range("L97").Value = 14040
$ws.Range("M97").Value = 37
$ws.Range("N97").Value = -15032

$ws.Range("H105").Value = 454001500
$ws.Range("J105").Value = 454001500
$ws.Range("L105").Value = 1362004500
$ws.Range("N105").Value = -1362009742

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H132").Value = 41667740
$ws.Range("I132").Value = 47619948
$ws.Range("J132").Value = 2301.6667
$ws.Range("K132").Value = 428579532
$ws.Range("L132").Value = 20715.0003
$ws.Range("M132").Value = -428577002
$ws.Range("N132").Value = -25775.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2270.1904
$ws.Range("I122").Value = 1231.2759
$ws.Range("J122").Value = 4587.769
$ws.Range("K122").Value = 3693.8277
$ws.Range("L122").Value = 13763.307
$ws.Range("M122").Value = -1243.8277
$ws.Range("N122").Value = -18663.307

$ws.Range("H126").Value = 3624.875
$ws.Range("I126").Value = 3571.4285
$ws.Range("J126").Value = 3999
$ws.Range("K126").Value = 10714.2855
$ws.Range("L126").Value = 11997
$ws.Range("M126").Value = -8244.2855
$ws.Range("N126").Value = -16937

$ws.Range("H132").Value = 2238.76
$ws.Range("I132").Value = 1883.55
$ws.Range("J132").Value = 3659.6
$ws.Range("K132").Value = 5650.65
$ws.Range("L132").Value = 10978.8
$ws.Range("M132").Value = -3120.65
$ws.Range("N132").Value = -16038.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1420.4375
$ws.Range("I61").Value = 1271.3334
$ws.Range("J61").Value = 1612.1428
$ws.Range("K61").Value = 1271.3334
$ws.Range("L61").Value = 1612.1428
$ws.Range("M61").Value = -1069.3334
$ws.Range("N61").Value = -2016.1428

$ws.Range("H113").Value = 1420.4375
$ws.Range("I113").Value = 1271.3334
$ws.Range("J113").Value = 1612.1428
$ws.Range("K113").Value = 1271.3334
$ws.Range("L113").Value = 1612.1428
$ws.Range("M113").Value = 898.6666
$ws.Range("N113").Value = -5952.1428

$ws.Range("H132").Value = 5949.484
$ws.Range("I132").Value = 7560.647
$ws.Range("J132").Value = 3993.0715
$ws.Range("K132").Value = 22681.941
$ws.Range("L132").Value = 11979.2145
$ws.Range("M132").Value = -20151.941
$ws.Range("N132").Value = -17039.2145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 14998
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 14998
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 14998
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -15296

$ws.Range("H42").Value = 11346.667
$ws.Range("I42").Value = 4044
$ws.Range("J42").Value = 14998
$ws.Range("K42").Value = 4044
$ws.Range("L42").Value = 14998
$ws.Range("M42").Value = -3666
$ws.Range("N42").Value = -15754

$ws.Range("H43").Value = 14850
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 14850
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 14850
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -15148

$ws.Range("H113").Value = 625.25
$ws.Range("J113").Value = 566.6667
$ws.Range("L113").Value = 1700.0001
$ws.Range("N113").Value = -6040.0001
